# Regenerate the localization handoff report.
#
# For the "Ready for handoff" files (rows 8, 9, 11, 12, 13, 14 on both the
# zh-cn and de-de sheets), a fresh handoff run:
#   - stamps the per-file "Priority" column with the handoff type ("ht")
#   - stamps the "Latest Handoff Datetime" column with the new generation time
# The Overview sheet's "Latest HO Xliff Generate Date" column shares that
# same handoff timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Sheets.Item("zh-cn")
$dede = $wb.Sheets.Item("de-de")
$overview = $wb.Sheets.Item("Overview")

$rows = @(8, 9, 11, 12, 13, 14)

$zhcnTimestamp = "2016-08-23 10:21:46"
$dedeTimestamp = "2016-08-23 10:21:51"

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = $zhcnTimestamp

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = $dedeTimestamp

    $overview.Range("G$r").Value = $dedeTimestamp
}
